$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"  = 101
    "E3"  = 43
    "E4"  = 46
    "E5"  = 140
    "E7"  = 35
    "F7"  = 19
    "H7"  = 24
    "E10" = 599
    "F10" = 302
    "H10" = 398
    "E11" = 385
    "F11" = 208
    "H11" = 272
    "E12" = 602
    "F12" = 330
    "H12" = 416
    "E14" = 130
    "F14" = 71
    "H14" = 105
    "E15" = 179
    "F16" = 105
    "H16" = 153
    "E17" = 105
    "F17" = 57
    "H17" = 81
    "E21" = 142
    "E22" = 176
    "F22" = 96
    "H22" = 138
    "E23" = 208
    "F23" = 102
    "H23" = 153
    "E24" = 226
    "F24" = 122
    "H24" = 152
    "E25" = 289
    "F25" = 146
    "H25" = 206
    "E26" = 163
    "F26" = 101
    "H26" = 126
    "E27" = 344
    "E28" = 207
    "F28" = 88
    "H28" = 140
    "E30" = 221
    "F30" = 132
    "H30" = 184
    "E31" = 75
    "F32" = 115
    "H32" = 153
    "E33" = 308
    "E34" = 227
    "F34" = 153
    "H34" = 192
    "E35" = 159
    "E36" = 79
    "F36" = 48
    "H36" = 58
    "E37" = 170
    "E38" = 96
    "E39" = 184
    "E40" = 274
    "F40" = 128
    "H40" = 208
    "E41" = 406
    "F41" = 197
    "H41" = 289
    "E42" = 399
    "F42" = 223
    "H42" = 284
    "E43" = 127
    "E44" = 324
    "F44" = 167
    "H44" = 235
    "E45" = 157
    "E46" = 339
    "F46" = 188
    "H46" = 251
    "E47" = 479
    "F47" = 251
    "H47" = 343
    "E48" = 229
    "E49" = 298
    "F49" = 135
    "H49" = 222
    "E50" = 249
    "F50" = 124
    "H50" = 195
    "E51" = 248
    "F51" = 112
    "H51" = 186
    "E52" = 30
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
